# Insert a new data row at row 26 (pushing existing rows 26-54 down to 27-55)
# and populate it with the new "Arveja Verde" price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("26").Insert()

$ws.Range("A26").Value = 11
$ws.Range("B26").Value = "Vega Monumental Concepción"
$ws.Range("C26").Value = "Bíobío"
$ws.Range("D26").Value = 45280
$ws.Range("E26").Value = 8
$ws.Range("F26").Value = 100112022
$ws.Range("G26").Value = "Arveja Verde"
$ws.Range("H26").Value = "Sin especificar"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 50
$ws.Range("K26").Value = 24000
$ws.Range("L26").Value = 24000
$ws.Range("M26").Value = 24000
$ws.Range("N26").Value = "`$/saco 25 kilos"
$ws.Range("O26").Value = "Región del Maule"
$ws.Range("P26").Value = 960
$ws.Range("Q26").Value = 25
$ws.Range("R26").Value = "Hortaliza"
